$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text, matching original inlineStr type,
# since assigning numeric-looking strings to .Value would auto-convert them to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.073.12"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.836.10"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "242.80"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "0.6279"
$ws.Range("E6").Value = "  -3.81%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.07609"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "22.60"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "0.07728"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.839.22"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "4.961"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "0.6656"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "0.00001015"
$ws.Range("E15").Value = "  +17.53%  "
$ws.Range("D16").Value = "82.88"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "6.063"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "29.012.60"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "226.86"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "7.200"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "158.67"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "8.505"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "0.1374"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "17.93"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "1.493"
$ws.Range("D29").Value = "4.101"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "4.019"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "1.191"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "0.05245"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "1.845"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "0.7346"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "2.697"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("D37").Value = "1.240.72"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "2.761"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "6.374"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "0.8973"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "102.13"
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("D45").Value = "1.981.48"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "64.16"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "0.5111"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "0.4043"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").Value = "8.853"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "0.05748"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "6.703"
$ws.Range("E51").Value = "  +0.31%  "

# Restore default number format/style so cell styling matches the original workbook.
$dRange.NumberFormat = "General"
$dRange.Style = "Normal"
